# ------------------------------------------------------------------
# Creacion workflow de modelos
# Renames the sole sheet "Modelos" -> "Ridge", adds five more model
# sheets (Lasso, Elastic net, Arboles, Forest, Boosting), removes the
# now-unused MAE_training column, adds a new "numero de modelo" column,
# fills in the first result row, bolds the header row and replicates
# the header row onto the new sheets.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the original sheet --------------------------------------
$ridge = $wb.Worksheets.Item(1)
$ridge.Name = "Ridge"

# --- 2. Create the other five sheets in the right creation order so the
#        internal sheetId sequence matches (Ridge=1, Elastic net=2,
#        Lasso=3, Arboles=4, Forest=5, Boosting=6) while re-fetching the
#        anchor sheet by name each time keeps the tab ORDER correct too:
#        Ridge, Lasso, Elastic net, Arboles, Forest, Boosting.
$elasticNet = $wb.Worksheets.Add($null, $wb.Worksheets.Item("Ridge"))
$elasticNet.Name = "Elastic net"

$lasso = $wb.Worksheets.Add($wb.Worksheets.Item("Elastic net"), $null)
$lasso.Name = "Lasso"

$arboles = $wb.Worksheets.Add($null, $wb.Worksheets.Item("Elastic net"))
$arboles.Name = "Arboles"

$forest = $wb.Worksheets.Add($null, $wb.Worksheets.Item("Arboles"))
$forest.Name = "Forest"

$boosting = $wb.Worksheets.Add($null, $wb.Worksheets.Item("Forest"))
$boosting.Name = "Boosting"

# --- 3. Fix up the "Ridge" sheet ----------------------------------------
# Drop the obsolete "MAE_training" column (column D); this shifts every
# column from E onward one place to the left.
$ridge.Columns.Item(4).Delete()

# New "numero de modelo" column at the end (column K).
$ridge.Range("K1").Value = "numero de modelo"

# Complete the first result row.
$ridge.Range("A2").Value = "tanda1_modelo1"
$ridge.Range("C2").Value = "test_1"
$ridge.Range("D2").Value = 30.245804
$ridge.Range("E2").Value = 4
$ridge.Range("K2").Value = 1

# Bold header row.
$ridge.Range("A1:K1").Font.Bold = $true

# Column widths; best-effort match of the authored widths (this engine
# quantizes ColumnWidth to 1/6-character steps on save, so these are the
# closest reachable values to the authored ones).
$ridge.Columns.Item(1).ColumnWidth = 16.333333333333332
$ridge.Columns.Item(3).ColumnWidth = 12.333333333333334
$ridge.Columns.Item(5).ColumnWidth = 17.5
$ridge.Columns.Item(8).ColumnWidth = 29.833333333333332
$ridge.Columns.Item(11).ColumnWidth = 18.333333333333332

# Restore the current selection shown in the saved file.
$ridge.Range("E34").Select()

# --- 4. Populate the header row (same 11 columns) on every new sheet ---
$headers = @("nombre resultados", "tanda", "datos_usados", "MAE_comp", "Leaderboard", "type", "mixture", "formula", "fecha", "submittor", "numero de modelo")

$newSheets = @($lasso, $elasticNet, $arboles, $forest, $boosting)
foreach ($sh in $newSheets) {
    for ($i = 0; $i -lt $headers.Count; $i++) {
        $sh.Cells.Item(1, $i + 1).Value = $headers[$i]
    }
    $sh.Range("A1:K1").Font.Bold = $true
}

# "Elastic net" has an extra, still-empty but explicitly formatted cell
# at E7 (work in progress row) in the authored file.
$elasticNet.Range("E7").Font.Bold = $false

# --- 5. Keep focus on the first (tab-selected) sheet --------------------
$ridge.Activate()
